$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D15: add student's note about light info (new shared string) ---
$ws.Range("D15").Value = "I read in light info from .py and set everything accordingly… just don’t have multiple light support"

# --- C18: "Minimap" feature earned (was 0, now full 0.11) ---
$ws.Range("C18").Value = 0.11
$ws.Range("C18").Font.Color = 255

# --- C46: "Game Music" feature - mark as earned style (value already 0.02) ---
$ws.Range("C46").Value = 0.02
$ws.Range("C46").Font.Color = 255

# --- Update active selection to match the new view state ---
$ws.Range("A14").Select()

Write-Host "Done applying edits"
